$d = $word.ActiveDocument

# 1. Insert the new answer text into the (currently empty) second list item,
#    right before the "_GoBack" bookmark that lives in that paragraph.
$apos = [char]0x2019
$answer2 = "It" + $apos + "s hard predict from the data. One hypothesis could be the loop counter register, since it needs to be written only once but could be read anywhere inside the loop. We can add more registers and use register renaming in order to remove a few long term dependencies which can be exploited to speed up out of order execution."

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Collapse(1)
$r2.InsertBefore($answer2)

# 2. Turn the third paragraph into a list item matching the other two
#    (same "List Paragraph" style + same numbering list, continuing it),
#    and drop the leading "3. " from its text.
$p1 = $d.Paragraphs.Item(1)
$p3 = $d.Paragraphs.Item(3)
$p3.Style = "List Paragraph"
$p3.Range.ListFormat.ApplyListTemplate($p1.Range.ListFormat.ListTemplate, $true)

$d.Content.Find.Execute("3. Architecture A", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Architecture A", 2)
